# Excel Download and Upload feature
# Update the D.O.B for row 2, clear the (now-unused) Full Name cell in row 3,
# append attachment note to the content/email body for row 3, and refresh the
# validation result message for row 3. Also move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update Date of Birth
$ws.Range("D2").Value = 35602

# Row 3: the name column is no longer populated
$ws.Range("B3").Clear()

# Row 3: email body now mentions the attachment
$ws.Range("F3").Value = "Hey Senior Dev Elango Hope you got this email.. With Attachment"

# Row 3: validation result reflects the invalid email / empty name case
$ws.Range("G3").Value = " Email ID is invalid | Name is empty |"

# Move the active cell selection to F12
$ws.Range("F12").Select()
